$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# DEJEUNER sheet: update the breakfast description cells (B2/C2) to mention
# "lait"/"milk" alongside the coffee beverage.
# ---------------------------------------------------------------------------
$wsDej = $wb.Worksheets.Item("DEJEUNER")
$wsDej.Range("B2").Value = "Œuf brouillé, bacon, jambon, saucisse, crêpe avec sirop, céréale froide, gruau, banane, orange et pomme, breuvage (jus purs 200 ml, lait, café)"
$wsDej.Range("C2").Value = "Scrambled egg, bacon, ham, sausage, pancake with syrup, cold cereal, oatmeal, banana, orange and apple, beverage (200 ml pure juice, milk, coffee)"

# ---------------------------------------------------------------------------
# QTY sheet: the "Jus pur 100 % ou lait / 100 % Pure Juice or milk" row is
# split into two separate rows - juice keeps the simplified label, and a new
# "Lait / Milk" row is inserted right after it with the same quantities.
# ---------------------------------------------------------------------------
$wsQty = $wb.Worksheets.Item("QTY")

$wsQty.Range("A2").Value = "Jus pur 100 %"
$wsQty.Range("B2").Value = "100 % Pure Juice"

$wsQty.Rows.Item(3).Insert()
$wsQty.Range("A3").Value = "Lait"
$wsQty.Range("B3").Value = "Milk"
$wsQty.Range("C3").Value = 1
$wsQty.Range("D3").Value = 1
$wsQty.Range("E3").Value = 1
$wsQty.Range("F3").Value = 1

# ---------------------------------------------------------------------------
# Active tab moves from SOUPER back to DEJEUNER.
# ---------------------------------------------------------------------------
$wsDej.Activate()
